# Commit: "changed â€™ to '"
#
# A handful of shared-string cells on the "Day1" sheet contain the
# mojibake sequence U+00E2 U+20AC U+2122 ("â€™") -- a right single quote
# that got double-encoded (UTF-8 bytes re-interpreted as Windows-1252/
# Latin-1). This script replaces every occurrence of that 3-character
# mojibake sequence with a plain straight apostrophe, leaving the other
# (still-broken) mojibake sequences such as â€œ / â€\x9d / ï¿½ untouched,
# exactly like the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Day1")

# Build the mojibake pattern from explicit code points so the .ps1 file
# itself stays plain ASCII (no encoding ambiguity when the file is read
# back by any tool).
$mojibake = [string]([char]0x00e2) + [string]([char]0x20ac) + [string]([char]0x2122)
$straightApostrophe = "'"

# Cells known (from the diff) to contain the mojibake sequence.
$targets = @("D2", "C9", "D9", "C11", "D13")

foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $old = $cell.Value
    if ($old -ne $null -and $old.Contains($mojibake)) {
        $new = $old.Replace($mojibake, $straightApostrophe)
        $cell.Value = $new
    }
}

Write-Host "Replaced mojibake a-euro-tm sequences with straight apostrophes in D2, C9, D9, C11, D13"
